$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize column B (drop auto bestFit, set an explicit custom width)
$ws.Columns.Item(2).ColumnWidth = 121

# Append the new Terms & Conditions rows (1557-1627)
$rows = @(
    @(1557, 'Welcome to Mentoryum!'),
    @(1558, 'These terms and conditions outline the rules and regulations for the use of Mentoryum''s Website,'),
    @(1559, 'By accessing this website, we assume you accept these terms and conditions. Do not continue to use Mentoryum if you do not agree to take all of the terms and conditions stated on this page.'),
    @(1560, '1. Terminology'),
    @(1561, 'The following terminology applies to these Terms and Conditions, Privacy Statement, and Disclaimer Notice, and all Agreements:'),
    @(1562, '"Client", "You" and "Your" refers to you, the person who logs on to this website and is compliant with the Company’s terms and conditions.'),
    @(1563, '"The Company", "Ourselves", "We", "Our", and "Us", refers to our Company.'),
    @(1564, '"Party", "Parties", or "Us", refers to both the Client and ourselves.'),
    @(1565, 'Any use of the above terminology or other words in the singular, plural, capitalization, and/or he/she or they, are taken as interchangeable and therefore as referring to the same.'),
    @(1566, '2. Cookies'),
    @(1567, 'We employ the use of cookies. By accessing Mentoryum, you agreed to use cookies in agreement with the Mentoryum''s Privacy Policy.'),
    @(1568, 'Most interactive websites use cookies to let us retrieve the user’s details for each visit. Cookies are used by our website to enable the functionality of certain areas to make it easier for people visiting our website. Some of our affiliate/advertising partners may also use cookies.'),
    @(1569, '3. License'),
    @(1570, 'Unless otherwise stated, Mentoryum and/or its licensors own the intellectual property rights for all material on Mentoryum. All intellectual property rights are reserved. You may access this from Mentoryum for your own personal use subjected to restrictions set in these terms and conditions.'),
    @(1571, 'You must not:'),
    @(1572, 'Republish material from Mentoryum'),
    @(1573, 'Sell, rent, or sub-license material from Mentoryum'),
    @(1574, 'Reproduce, duplicate, or copy material from Mentoryum'),
    @(1575, 'Redistribute content from Mentoryum'),
    @(1576, 'This Agreement shall begin on the date hereof.'),
    @(1577, '4. User Comments'),
    @(1578, 'Parts of this website offer an opportunity for users to post and exchange opinions and information in certain areas of the website. Mentoryum does not filter, edit, publish or review Comments prior to their presence on the website. Comments do not reflect the views and opinions of Mentoryum, its agents, and/or affiliates. Comments reflect the views and opinions of the person who posts their views and opinions.'),
    @(1579, 'Mentoryum shall not be liable for the Comments or any liability, damages, or expenses caused and/or suffered as a result of any use of and/or posting of and/or appearance of the Comments on this website.'),
    @(1580, 'Mentoryum reserves the right to monitor all Comments and to remove any Comments that can be considered inappropriate, offensive, or causes a breach of these Terms and Conditions.'),
    @(1581, 'You warrant and represent that:'),
    @(1582, 'You are entitled to post the Comments on our website and have all necessary licenses and consents to do so;'),
    @(1583, 'The Comments do not invade any intellectual property right, including without limitation copyright, patent, or trademark of any third party;'),
    @(1584, 'The Comments do not contain any defamatory, libelous, offensive, indecent, or otherwise unlawful material which is an invasion of privacy;'),
    @(1585, 'The Comments will not be used to solicit or promote business or custom or present commercial activities or unlawful activity.'),
    @(1586, 'You hereby grant Mentoryum a non-exclusive license to use, reproduce, edit and authorize others to use, reproduce and edit any of your Comments in any and all forms, formats, or media.'),
    @(1587, '5. Hyperlinking to our Content'),
    @(1588, 'The following organizations may link to our Website without prior written approval:'),
    @(1589, 'Government agencies;'),
    @(1590, 'Search engines;'),
    @(1591, 'News organizations;'),
    @(1592, 'Online directory distributors may link to our Website in the same manner as they hyperlink to the Websites of other listed businesses; and'),
    @(1593, 'System-wide Accredited Businesses except soliciting non-profit organizations, charity shopping malls, and charity fundraising groups which may not hyperlink to our Web site.'),
    @(1594, 'These organizations may link to our home page, to publications, or to other Website information so long as the link: (a) is not in any way deceptive; (b) does not falsely imply sponsorship, endorsement, or approval of the linking party and its products and/or services; and (c) fits within the context of the linking party’s site.'),
    @(1595, 'We may consider and approve other link requests from the following types of organizations:'),
    @(1596, 'commonly-known consumer and/or business information sources;'),
    @(1597, 'dot.com community sites;'),
    @(1598, 'associations or other groups representing charities;'),
    @(1599, 'online directory distributors;'),
    @(1600, 'internet portals;'),
    @(1601, 'accounting, law, and consulting firms; and'),
    @(1602, 'educational institutions and trade associations.'),
    @(1603, 'We will approve link requests from these organizations if we decide that: (a) the link would not make us look unfavorably to ourselves or to our accredited businesses; (b) the organization does not have any negative records with us; (c) the benefit to us from the visibility of the hyperlink compensates the absence of Mentoryum; and (d) the link is in the context of general resource information.'),
    @(1604, 'These organizations may link to our home page so long as the link: (a) is not in any way deceptive; (b) does not falsely imply sponsorship, endorsement, or approval of the linking party and its products or services; and (c) fits within the context of the linking party’s site.'),
    @(1605, 'If you are one of the organizations listed in paragraph 2 above and are interested in linking to our website, you must inform us by sending an e-mail to Mentoryum. Please include your name, your organization name, contact information as well as the URL of your site, a list of any URLs from which you intend to link to our Website, and a list of the URLs on our site to which you would like to link. Wait 2-3 weeks for a response.'),
    @(1606, 'Approved organizations may hyperlink to our Website as follows:'),
    @(1607, 'By use of our corporate name; or'),
    @(1608, 'By use of the uniform resource locator being linked to; or'),
    @(1609, 'By use of any other description of our Website being linked to that makes sense within the context and format of content on the linking party’s site.'),
    @(1610, 'No use of Mentoryum''s logo or other artwork will be allowed for linking absent a trademark license agreement.'),
    @(1611, '6. iFrames'),
    @(1612, 'Without prior approval and written permission, you may not create frames around our Webpages that alter in any way the visual presentation or appearance of our Website.'),
    @(1613, '7. Content Liability'),
    @(1614, 'We shall not be held responsible for any content that appears on your Website. You agree to protect and defend us against all claims that are rising on your Website. No link(s) should appear on any Website that may be interpreted as libelous, obscene, or criminal, or which infringes, otherwise violates, or advocates the infringement or other violation of, any third party rights.'),
    @(1615, '8. Reservation of Rights'),
    @(1616, 'We reserve the right to request that you remove all links or any particular link to our Website. You approve to immediately remove all links to our Website upon request. We also reserve the right to amend these terms and conditions and it’s linking policy at any time. By continuously linking to our Website, you agree to be bound to and follow these linking terms and conditions.'),
    @(1617, '9. Removal of links from our website'),
    @(1618, 'If you find any link on our Website that is offensive for any reason, you are free to contact and inform us at any moment. We will consider requests to remove links but we are not obligated to or so or to respond to you directly.'),
    @(1619, 'We do not ensure that the information on this website is correct, we do not warrant its completeness or accuracy; nor do we promise to ensure that the website remains available or that the material on the website is kept up to date.'),
    @(1620, '10. Disclaimer'),
    @(1621, 'To the maximum extent permitted by applicable law, we exclude all representations, warranties, and conditions relating to our website and the use of this website. Nothing in this disclaimer will:'),
    @(1622, 'limit or exclude our or your liability for death or personal injury;'),
    @(1623, 'limit or exclude our or your liability for fraud or fraudulent misrepresentation;'),
    @(1624, 'limit any of our or your liabilities in any way that is not permitted under applicable law; or'),
    @(1625, 'exclude any of our or your liabilities that may not be excluded under applicable law.'),
    @(1626, 'The limitations and prohibitions of liability set in this Section and elsewhere in this disclaimer: (a) are subject to the preceding paragraph; and (b) govern all liabilities arising under the disclaimer, including liabilities arising in contract, in tort, and for breach of statutory duty.'),
    @(1627, 'As long as the website and the information and services on the website are provided free of charge, we will not be liable for any loss or damage of any nature.')
)

foreach ($item in $rows) {
    $r = $item[0]
    $t = $item[1]
    $ws.Cells.Item($r, 1).Value2 = ($r - 1)
    $ws.Cells.Item($r, 2).Value2 = $t
    $ws.Cells.Item($r, 3).Value2 = 1
}

# "10. Disclaimer" heading gets center vertical alignment (matches the other section headings style)
$ws.Range("B1551").Copy()
$ws.Cells.Item(1620, 2).PasteSpecial(-4122)
$ws.Cells.Item(1620, 2).Value2 = "10. Disclaimer"
$ws.Cells.Item(1620, 2).Font.Name = "Calibri"
$ws.Cells.Item(1620, 2).VerticalAlignment = -4108

# Restore final selection/scroll state
$ws.Cells.Item(1596, 2).Select()

Write-Host "done"